$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 44.30160133333333
$ws.Cells.Item(2, 8).Value = 132.904804
$ws.Cells.Item(2, 9).Value = 0.1310981746002036
$ws.Cells.Item(2, 10).Value = 0.1395903267786693
$ws.Cells.Item(2, 13).Value = 15.70818033333333
$ws.Cells.Item(2, 14).Value = 47.12454099999999
$ws.Cells.Item(2, 15).Value = 0.3220467100482788
$ws.Cells.Item(2, 16).Value = 0.334408980496766
$ws.Cells.Item(2, 17).Value = 695.8975427994403
$ws.Cells.Item(2, 18).Value = 6263.077885194963
$ws.Cells.Item(2, 19).Value = 0.04221973582333041
$ws.Cells.Item(2, 20).Value = 0.04668025886526523
$ws.Cells.Item(3, 7).Value = 44.30160133333333
$ws.Cells.Item(3, 8).Value = 132.904804
$ws.Cells.Item(3, 9).Value = 0.1310981746002036
$ws.Cells.Item(3, 10).Value = 0.1395903267786693
$ws.Cells.Item(3, 15).Value = 0.1200026410479322
$ws.Cells.Item(3, 16).Value = 0.1246091315254933
$ws.Cells.Item(3, 17).Value = 259.3087910203466
$ws.Cells.Item(3, 18).Value = 2333.77911918312
$ws.Cells.Item(3, 19).Value = 0.01573212718858738
$ws.Cells.Item(3, 20).Value = 0.0173942293892498
$ws.Cells.Item(4, 7).Value = 44.30160133333333
$ws.Cells.Item(4, 8).Value = 132.904804
$ws.Cells.Item(4, 9).Value = 0.1310981746002036
$ws.Cells.Item(4, 10).Value = 0.1395903267786693
$ws.Cells.Item(4, 13).Value = 10.959131
$ws.Cells.Item(4, 14).Value = 32.877393
$ws.Cells.Item(4, 15).Value = 0.2246824271585863
$ws.Cells.Item(4, 16).Value = 0.2333072161810874
$ws.Cells.Item(4, 17).Value = 485.5070525217747
$ws.Cells.Item(4, 18).Value = 4369.563472695972
$ws.Cells.Item(4, 19).Value = 0.02945545606523388
$ws.Cells.Item(4, 20).Value = 0.03256743054653964
$ws.Cells.Item(5, 7).Value = 44.30160133333333
$ws.Cells.Item(5, 8).Value = 132.904804
$ws.Cells.Item(5, 9).Value = 0.1310981746002036
$ws.Cells.Item(5, 10).Value = 0.1395903267786693
$ws.Cells.Item(5, 13).Value = 5.4093935
$ws.Cells.Item(5, 14).Value = 10.818787
$ws.Cells.Item(5, 15).Value = 0.1109025579706895
$ws.Cells.Item(5, 16).Value = 0.07677315161290731
$ws.Cells.Item(5, 17).Value = 239.6447942921247
$ws.Cells.Item(5, 18).Value = 1437.868765752748
$ws.Cells.Item(5, 19).Value = 0.01453912290845065
$ws.Cells.Item(5, 20).Value = 0.01071678932147406
$ws.Cells.Item(6, 7).Value = 44.30160133333333
$ws.Cells.Item(6, 8).Value = 132.904804
$ws.Cells.Item(6, 9).Value = 0.1310981746002036
$ws.Cells.Item(6, 10).Value = 0.1395903267786693
$ws.Cells.Item(6, 13).Value = 10.84612833333333
$ws.Cells.Item(6, 14).Value = 32.538385
$ws.Cells.Item(6, 15).Value = 0.2223656637745133
$ws.Cells.Item(6, 16).Value = 0.230901520183746
$ws.Cells.Item(6, 17).Value = 480.5008534335044
$ws.Cells.Item(6, 18).Value = 4324.507680901539
$ws.Cells.Item(6, 19).Value = 0.02915173261460131
$ws.Cells.Item(6, 20).Value = 0.03223161865614062
$ws.Cells.Item(7, 9).Value = 0.2142454163706631
$ws.Cells.Item(7, 10).Value = 0.2281236010586413
$ws.Cells.Item(7, 13).Value = 15.70818033333333
$ws.Cells.Item(7, 14).Value = 47.12454099999999
$ws.Cells.Item(7, 15).Value = 0.3220467100482788
$ws.Cells.Item(7, 16).Value = 0.334408980496766
$ws.Cells.Item(7, 17).Value = 1137.261134741657
$ws.Cells.Item(7, 18).Value = 10235.35021267492
$ws.Cells.Item(7, 19).Value = 0.06899703148509571
$ws.Cells.Item(7, 20).Value = 0.07628658085727121
$ws.Cells.Item(8, 9).Value = 0.2142454163706631
$ws.Cells.Item(8, 10).Value = 0.2281236010586413
$ws.Cells.Item(8, 15).Value = 0.1200026410479322
$ws.Cells.Item(8, 16).Value = 0.1246091315254933
$ws.Cells.Item(8, 19).Value = 0.02571001579689347
$ws.Cells.Item(8, 20).Value = 0.0284262838083854
$ws.Cells.Item(9, 9).Value = 0.2142454163706631
$ws.Cells.Item(9, 10).Value = 0.2281236010586413
$ws.Cells.Item(9, 13).Value = 10.959131
$ws.Cells.Item(9, 14).Value = 32.877393
$ws.Cells.Item(9, 15).Value = 0.2246824271585863
$ws.Cells.Item(9, 16).Value = 0.2333072161810874
$ws.Cells.Item(9, 17).Value = 793.4333253352521
$ws.Cells.Item(9, 18).Value = 7140.899928017269
$ws.Cells.Item(9, 19).Value = 0.04813718015776251
$ws.Cells.Item(9, 20).Value = 0.05322288230819655
$ws.Cells.Item(10, 9).Value = 0.2142454163706631
$ws.Cells.Item(10, 10).Value = 0.2281236010586413
$ws.Cells.Item(10, 13).Value = 5.4093935
$ws.Cells.Item(10, 14).Value = 10.818787
$ws.Cells.Item(10, 15).Value = 0.1109025579706895
$ws.Cells.Item(10, 16).Value = 0.07677315161290731
$ws.Cells.Item(10, 17).Value = 391.636259549402
$ws.Cells.Item(10, 18).Value = 2349.817557296412
$ws.Cells.Item(10, 19).Value = 0.02376036470900197
$ws.Cells.Item(10, 20).Value = 0.01751376781055745
$ws.Cells.Item(11, 9).Value = 0.2142454163706631
$ws.Cells.Item(11, 10).Value = 0.2281236010586413
$ws.Cells.Item(11, 13).Value = 10.84612833333333
$ws.Cells.Item(11, 14).Value = 32.538385
$ws.Cells.Item(11, 15).Value = 0.2223656637745133
$ws.Cells.Item(11, 16).Value = 0.230901520183746
$ws.Cells.Item(11, 17).Value = 785.2520122744734
$ws.Cells.Item(11, 18).Value = 7067.26811047026
$ws.Cells.Item(11, 19).Value = 0.04764082422190948
$ws.Cells.Item(11, 20).Value = 0.05267408627423067
$ws.Cells.Item(12, 7).Value = 82.35175066666666
$ws.Cells.Item(12, 8).Value = 247.055252
$ws.Cells.Item(12, 9).Value = 0.2436969288378267
$ws.Cells.Item(12, 10).Value = 0.2594828954344383
$ws.Cells.Item(12, 13).Value = 15.70818033333333
$ws.Cells.Item(12, 14).Value = 47.12454099999999
$ws.Cells.Item(12, 15).Value = 0.3220467100482788
$ws.Cells.Item(12, 16).Value = 0.334408980496766
$ws.Cells.Item(12, 17).Value = 1293.596150237703
$ws.Cells.Item(12, 18).Value = 11642.36535213933
$ws.Cells.Item(12, 19).Value = 0.07848179418109159
$ws.Cells.Item(12, 20).Value = 0.08677341051857944
$ws.Cells.Item(13, 7).Value = 82.35175066666666
$ws.Cells.Item(13, 8).Value = 247.055252
$ws.Cells.Item(13, 9).Value = 0.2436969288378267
$ws.Cells.Item(13, 10).Value = 0.2594828954344383
$ws.Cells.Item(13, 15).Value = 0.1200026410479322
$ws.Cells.Item(13, 16).Value = 0.1246091315254933
$ws.Cells.Item(13, 17).Value = 482.0262081071733
$ws.Cells.Item(13, 18).Value = 4338.23587296456
$ws.Cells.Item(13, 19).Value = 0.0292442750758092
$ws.Cells.Item(13, 20).Value = 0.03233393824580575
$ws.Cells.Item(14, 7).Value = 82.35175066666666
$ws.Cells.Item(14, 8).Value = 247.055252
$ws.Cells.Item(14, 9).Value = 0.2436969288378267
$ws.Cells.Item(14, 10).Value = 0.2594828954344383
$ws.Cells.Item(14, 13).Value = 10.959131
$ws.Cells.Item(14, 14).Value = 32.877393
$ws.Cells.Item(14, 15).Value = 0.2246824271585863
$ws.Cells.Item(14, 16).Value = 0.2333072161810874
$ws.Cells.Item(14, 17).Value = 902.5036236353374
$ws.Cells.Item(14, 18).Value = 8122.532612718037
$ws.Cells.Item(14, 19).Value = 0.05475441746237618
$ws.Cells.Item(14, 20).Value = 0.06053923198041698
$ws.Cells.Item(15, 7).Value = 82.35175066666666
$ws.Cells.Item(15, 8).Value = 247.055252
$ws.Cells.Item(15, 9).Value = 0.2436969288378267
$ws.Cells.Item(15, 10).Value = 0.2594828954344383
$ws.Cells.Item(15, 13).Value = 5.4093935
$ws.Cells.Item(15, 14).Value = 10.818787
$ws.Cells.Item(15, 15).Value = 0.1109025579706895
$ws.Cells.Item(15, 16).Value = 0.07677315161290731
$ws.Cells.Item(15, 17).Value = 445.4730247698873
$ws.Cells.Item(15, 18).Value = 2672.838148619324
$ws.Cells.Item(15, 19).Value = 0.02702661277771606
$ws.Cells.Item(15, 20).Value = 0.0199213196721443
$ws.Cells.Item(16, 7).Value = 82.35175066666666
$ws.Cells.Item(16, 8).Value = 247.055252
$ws.Cells.Item(16, 9).Value = 0.2436969288378267
$ws.Cells.Item(16, 10).Value = 0.2594828954344383
$ws.Cells.Item(16, 13).Value = 10.84612833333333
$ws.Cells.Item(16, 14).Value = 32.538385
$ws.Cells.Item(16, 15).Value = 0.2223656637745133
$ws.Cells.Item(16, 16).Value = 0.230901520183746
$ws.Cells.Item(16, 17).Value = 893.1976562053354
$ws.Cells.Item(16, 18).Value = 8038.778905848019
$ws.Cells.Item(16, 19).Value = 0.05418982934083365
$ws.Cells.Item(16, 20).Value = 0.05991499501749179
$ws.Cells.Item(17, 7).Value = 61.6746195
$ws.Cells.Item(17, 8).Value = 123.349239
$ws.Cells.Item(17, 9).Value = 0.1825087534596294
$ws.Cells.Item(17, 10).Value = 0.1295540872992837
$ws.Cells.Item(17, 13).Value = 15.70818033333333
$ws.Cells.Item(17, 14).Value = 47.12454099999999
$ws.Cells.Item(17, 15).Value = 0.3220467100482788
$ws.Cells.Item(17, 16).Value = 0.334408980496766
$ws.Cells.Item(17, 17).Value = 968.7960450957163
$ws.Cells.Item(17, 18).Value = 5812.776270574298
$ws.Cells.Item(17, 19).Value = 0.05877634360668607
$ws.Cells.Item(17, 20).Value = 0.04332405025294248
$ws.Cells.Item(18, 7).Value = 61.6746195
$ws.Cells.Item(18, 8).Value = 123.349239
$ws.Cells.Item(18, 9).Value = 0.1825087534596294
$ws.Cells.Item(18, 10).Value = 0.1295540872992837
$ws.Cells.Item(18, 15).Value = 0.1200026410479322
$ws.Cells.Item(18, 16).Value = 0.1246091315254933
$ws.Cells.Item(18, 17).Value = 360.99758333457
$ws.Cells.Item(18, 18).Value = 2165.98550000742
$ws.Cells.Item(18, 19).Value = 0.02190153242952147
$ws.Cells.Item(18, 20).Value = 0.01614362230394169
$ws.Cells.Item(19, 7).Value = 61.6746195
$ws.Cells.Item(19, 8).Value = 123.349239
$ws.Cells.Item(19, 9).Value = 0.1825087534596294
$ws.Cells.Item(19, 10).Value = 0.1295540872992837
$ws.Cells.Item(19, 13).Value = 10.959131
$ws.Cells.Item(19, 14).Value = 32.877393
$ws.Cells.Item(19, 15).Value = 0.2246824271585863
$ws.Cells.Item(19, 16).Value = 0.2333072161810874
$ws.Cells.Item(19, 17).Value = 675.9002344756545
$ws.Cells.Item(19, 18).Value = 4055.401406853928
$ws.Cells.Item(19, 19).Value = 0.04100650970499757
$ws.Cells.Item(19, 20).Value = 0.03022590345267745
$ws.Cells.Item(20, 7).Value = 61.6746195
$ws.Cells.Item(20, 8).Value = 123.349239
$ws.Cells.Item(20, 9).Value = 0.1825087534596294
$ws.Cells.Item(20, 10).Value = 0.1295540872992837
$ws.Cells.Item(20, 13).Value = 5.4093935
$ws.Cells.Item(20, 14).Value = 10.818787
$ws.Cells.Item(20, 15).Value = 0.1109025579706895
$ws.Cells.Item(20, 16).Value = 0.07677315161290731
$ws.Cells.Item(20, 17).Value = 333.6222858382732
$ws.Cells.Item(20, 18).Value = 1334.489143353093
$ws.Cells.Item(20, 19).Value = 0.02024068761071482
$ws.Cells.Item(20, 20).Value = 0.009946275586299735
$ws.Cells.Item(21, 7).Value = 61.6746195
$ws.Cells.Item(21, 8).Value = 123.349239
$ws.Cells.Item(21, 9).Value = 0.1825087534596294
$ws.Cells.Item(21, 10).Value = 0.1295540872992837
$ws.Cells.Item(21, 13).Value = 10.84612833333333
$ws.Cells.Item(21, 14).Value = 32.538385
$ws.Cells.Item(21, 15).Value = 0.2223656637745133
$ws.Cells.Item(21, 16).Value = 0.230901520183746
$ws.Cells.Item(21, 17).Value = 668.9308380065024
$ws.Cells.Item(21, 18).Value = 4013.585028039015
$ws.Cells.Item(21, 19).Value = 0.04058368010770948
$ws.Cells.Item(21, 20).Value = 0.02991423570342234
$ws.Cells.Item(22, 7).Value = 77.19964866666666
$ws.Cells.Item(22, 8).Value = 231.598946
$ws.Cells.Item(22, 9).Value = 0.2284507267316773
$ws.Cells.Item(22, 10).Value = 0.2432490894289675
$ws.Cells.Item(22, 13).Value = 15.70818033333333
$ws.Cells.Item(22, 14).Value = 47.12454099999999
$ws.Cells.Item(22, 15).Value = 0.3220467100482788
$ws.Cells.Item(22, 16).Value = 0.334408980496766
$ws.Cells.Item(22, 17).Value = 1212.666002925976
$ws.Cells.Item(22, 18).Value = 10913.99402633378
$ws.Cells.Item(22, 19).Value = 0.07357180495207503
$ws.Cells.Item(22, 20).Value = 0.0813446800027077
$ws.Cells.Item(23, 7).Value = 77.19964866666666
$ws.Cells.Item(23, 8).Value = 231.598946
$ws.Cells.Item(23, 9).Value = 0.2284507267316773
$ws.Cells.Item(23, 10).Value = 0.2432490894289675
$ws.Cells.Item(23, 15).Value = 0.1200026410479322
$ws.Cells.Item(23, 16).Value = 0.1246091315254933
$ws.Cells.Item(23, 17).Value = 451.8696155546533
$ws.Cells.Item(23, 18).Value = 4066.82653999188
$ws.Cells.Item(23, 19).Value = 0.02741469055712072
$ws.Cells.Item(23, 20).Value = 0.0303110577781107
$ws.Cells.Item(24, 7).Value = 77.19964866666666
$ws.Cells.Item(24, 8).Value = 231.598946
$ws.Cells.Item(24, 9).Value = 0.2284507267316773
$ws.Cells.Item(24, 10).Value = 0.2432490894289675
$ws.Cells.Item(24, 13).Value = 10.959131
$ws.Cells.Item(24, 14).Value = 32.877393
$ws.Cells.Item(24, 15).Value = 0.2246824271585863
$ws.Cells.Item(24, 16).Value = 0.2333072161810874
$ws.Cells.Item(24, 17).Value = 846.0410628919753
$ws.Cells.Item(24, 18).Value = 7614.369566027779
$ws.Cells.Item(24, 19).Value = 0.05132886376821618
$ws.Cells.Item(24, 20).Value = 0.05675176789325678
$ws.Cells.Item(25, 7).Value = 77.19964866666666
$ws.Cells.Item(25, 8).Value = 231.598946
$ws.Cells.Item(25, 9).Value = 0.2284507267316773
$ws.Cells.Item(25, 10).Value = 0.2432490894289675
$ws.Cells.Item(25, 13).Value = 5.4093935
$ws.Cells.Item(25, 14).Value = 10.818787
$ws.Cells.Item(25, 15).Value = 0.1109025579706895
$ws.Cells.Item(25, 16).Value = 0.07677315161290731
$ws.Cells.Item(25, 17).Value = 417.6032776997503
$ws.Cells.Item(25, 18).Value = 2505.619666198502
$ws.Cells.Item(25, 19).Value = 0.02533576996480598
$ws.Cells.Item(25, 20).Value = 0.01867499922243177
$ws.Cells.Item(26, 7).Value = 77.19964866666666
$ws.Cells.Item(26, 8).Value = 231.598946
$ws.Cells.Item(26, 9).Value = 0.2284507267316773
$ws.Cells.Item(26, 10).Value = 0.2432490894289675
$ws.Cells.Item(26, 13).Value = 10.84612833333333
$ws.Cells.Item(26, 14).Value = 32.538385
$ws.Cells.Item(26, 15).Value = 0.2223656637745133
$ws.Cells.Item(26, 16).Value = 0.230901520183746
$ws.Cells.Item(26, 17).Value = 837.3172967269121
$ws.Cells.Item(26, 18).Value = 7535.855670542209
$ws.Cells.Item(26, 19).Value = 0.05079959748945936
$ws.Cells.Item(26, 20).Value = 0.05616658453246057
